$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V2:V3").Value = 0.05
$ws.Range("X2:X3").Value = 0.1032696664629665
$ws.Range("AB2:AB3").Value = 0.1032696664629665
$ws.Range("AD2:AD3").Value = 0
$ws.Range("AF2:AF3").Value = 0
$ws.Range("AG2:AG3").Value = -0.266
$ws.Range("AH2:AH3").Value = 0
$ws.Range("AI2:AI3").Value = 0
$ws.Range("AJ2:AJ3").Value = -0.05263157894736842
$ws.Range("AK2:AK3").Value = -0.02902662592754256
$ws.Range("AN2:AN3").Value = 0
$ws.Range("AP2:AP3").Value = -0.3384223918575064
